# Update "想去人数" (attendee interest counts) figures in the 展览 and
# 全部类型 sheets to match the newly scraped output.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 10092
$ws1.Range("F3").Value = 225
$ws1.Range("F4").Value = 49
$ws1.Range("F5").Value = 619
$ws1.Range("F6").Value = 480

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 10092
$ws4.Range("F3").Value = 225
$ws4.Range("F4").Value = 49
$ws4.Range("F5").Value = 619
$ws4.Range("F7").Value = 480
